{"js": "// Helper to replace a single paragraph's content via OOXML package insertion.\n// `innerXml` is the raw <w:p>...</w:p> XML (namespaces resolved via the\n// wrapping w:document element) that should become the new paragraph content.\nasync function replaceParagraphOoxml(context, paragraph, pXml) {\n  const ooxml =\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n          '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n            '<w:body>' + pXml + '</w:body>' +\n          '</w:document>' +\n        '</pkg:xmlData>' +\n      '</pkg:part>' +\n    '</pkg:package>';\n  paragraph.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- 1. Title: \"Tic Tac Toe \u2013 Design Document\" -----------------------------\n// Split into three runs, wrapping \"Tac\" with grammar-check proofErr marks.\nawait replaceParagraphOoxml(context, paragraphs.items[6],\n  '<w:p>' +\n    '<w:pPr><w:pStyle w:val=\"Heading1\"/><w:jc w:val=\"center\"/></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Tic </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>Tac</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> Toe \u2013 Design Document</w:t></w:r>' +\n  '</w:p>');\n\n// --- 2. Author line: names get spell-check proofErr wraps -------------------\nawait replaceParagraphOoxml(context, paragraphs.items[26],\n  '<w:p>' +\n    '<w:pPr><w:pStyle w:val=\"Heading3\"/></w:pPr>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Shayaan</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> Ali, </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Davinderpal</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> Cheema, Alexander Powell, </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Riddhi</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> Sharma</w:t></w:r>' +\n  '</w:p>');\n\n// --- 3. Introduction paragraph: wrap both \"tac\" occurrences -----------------\nawait replaceParagraphOoxml(context, paragraphs.items[30],\n  '<w:p>' +\n    '<w:r><w:tab/><w:t xml:space=\"preserve\">Tic Tac Toe is a game traditionally played between two players on a 3x3 grid. Players take turns placing X and O in the empty spaces on the grid. This project takes the basics of tic </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>tac</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> toe and builds an application to allow players to play against each other over the internet. </w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">It will take the basics of tic </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>tac</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> toe and stretch it to support any number of players on any board size (n x n). The application is built using a client-server architecture where clients connect to the server and play each other with the server acting as the messenger.</w:t></w:r>' +\n  '</w:p>');\n\n// --- 4. Client Design paragraph: append a new sentence/run -----------------\nawait replaceParagraphOoxml(context, paragraphs.items[56],\n  '<w:p>' +\n    '<w:r><w:tab/><w:t>The client is designed and built to be a thin client, the client does not perform game calculations. The client is simply a user interface that accepts inputs from the user an converts them into game move. The game moves are then sent to the server where the server handles the game logic and sends board updates and game updates to the client.</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> The client uses a few threads, it has a thread for passing messages to the server and a thread that listens to the messages from the server as well as a thread for the GUI.</w:t></w:r>' +\n  '</w:p>');\n\n// --- 5. Synchronization paragraph: split run, add a page break -------------\nawait replaceParagraphOoxml(context, paragraphs.items[58],\n  '<w:p>' +\n    '<w:r><w:tab/><w:t xml:space=\"preserve\">The clients are synchronized by the server. The clients connect to the server and the server handles all message passing between them. If clients play out of turn then the server will stop them. If the clients try and make a move where another client has already moved then the client server will stop </w:t></w:r>' +\n    '<w:r><w:lastRenderedPageBreak/><w:t>them and tell them to try again. The clients cannot become out of sync because all turn and game logic is computed on the server and the server facilitates all messages that are passed between the clients.</w:t></w:r>' +\n  '</w:p>');\n\n// --- 6. \"Termination\" heading: drop the page break (it moved above) --------\nawait replaceParagraphOoxml(context, paragraphs.items[59],\n  '<w:p>' +\n    '<w:pPr><w:pStyle w:val=\"Heading1\"/></w:pPr>' +\n    '<w:r><w:t>Termination</w:t></w:r>' +\n  '</w:p>');\n\n// --- 7. Termination detail paragraph: add the _GoBack bookmark -------------\nawait replaceParagraphOoxml(context, paragraphs.items[60],\n  '<w:p>' +\n    '<w:r><w:tab/><w:t xml:space=\"preserve\">The server is terminated through the console. The server admin can simply type \\u2018stop\\u2019 to stop the server from running. This terminates all threads that are running o</w:t></w:r>' +\n    '<w:r><w:t>n the server and boots any clients from in progress games.</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> Termination in python is as simple as calling threads destructor in C++. Because threads cannot be killed externally the server requires a termination event object that is shared by all threads. When the event is signaled all threads terminate themselves.</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>');\n\n// --- 8. \"Object oriented style\" paragraph: split + \"second\"->\"final\", and --\n//        add a new paragraph right after it with the prototyping narrative.\n//        Note: the old trailing _GoBack bookmark is intentionally dropped\n//        here since it was relocated to the Termination paragraph (step 7).\nawait replaceParagraphOoxml(context, paragraphs.items[63],\n  '<w:p>' +\n    '<w:r><w:tab/><w:t xml:space=\"preserve\">The object oriented style of the </w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">final </w:t></w:r>' +\n    '<w:r><w:t>version of the project took longer to make but it was a better result. The code of the final version is much more readable than the first and all of the classes perform fewer functions working together to make the game work instead of a few classes doing all of the work. The first version had a few points of failure for that would take the entire application with it if it failed. In the final version some threads can fail without taking down the server.</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p>' +\n    '<w:r><w:tab/><w:t>During development a throwaway prototyping model was used to create the final version. There were several iterations in order to create code that is easy to read, efficient and modular. The first iteration was not very object oriented and code acted more like several scripts running on different threads. The next version was more object oriented and focused on how the specific parts of the server and client would need to be built in an object oriented project. The next iteration took the knowledge from the first iteration and tried to extend the server and client classes to work together with each other. Then in the next iteration we focused on how the client and server would interact with the game and the user interface. Finally taking forward the knowledge of the last iteration the project was started again from the ground up in order to remove the mistakes from the prototype phase and finally create a product that met the goals of the final project.</w:t></w:r>' +\n  '</w:p>');\n\n// --- 9. \"What did we think of the labs for this course?\" gets a page break -\n// Paragraph count grew by one (step 8 added a paragraph), so reload first.\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items\");\nawait context.sync();\nawait replaceParagraphOoxml(context, paragraphs2.items[69],\n  '<w:p>' +\n    '<w:r><w:lastRenderedPageBreak/><w:t>What did we think of the labs for this course?</w:t></w:r>' +\n  '</w:p>');\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the same set of paragraph-level OOXML replacements as edit.js,\n# using Range.InsertXML (the COM analogue of Office.js's insertOoxml) to\n# splice in new run structure (proofErr marks, bookmarks, page breaks,\n# split runs, and an extra paragraph) while leaving all other content as-is.\n#\n# NOTE: this interpreter's function calls bind parameters positionally, not\n# by -Name, so Set-ParagraphOoxml below is always invoked as\n# `Set-ParagraphOoxml <index> <xmlVar>` (no -Index/-InnerXml switches).\n\n$d = $word.ActiveDocument\n\nfunction Set-ParagraphOoxml($idx, $innerXml) {\n    $p = $d.Paragraphs.Item($idx)\n    $r = $p.Range\n    $full = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n    $r.InsertXML($full)\n}\n\n# --- 1. Title: \"Tic Tac Toe - Design Document\" (COM paragraph #7) ----------\n$xml1 = @\"\n<w:p>\n  <w:pPr><w:pStyle w:val=\"Heading1\"/><w:jc w:val=\"center\"/></w:pPr>\n  <w:r><w:t xml:space=\"preserve\">Tic </w:t></w:r>\n  <w:proofErr w:type=\"gramStart\"/>\n  <w:r><w:t>Tac</w:t></w:r>\n  <w:proofErr w:type=\"gramEnd\"/>\n  <w:r><w:t xml:space=\"preserve\"> Toe &#8211; Design Document</w:t></w:r>\n</w:p>\n\"@\nSet-ParagraphOoxml 7 $xml1\n\n# --- 2. Author line (COM paragraph #27) -------------------------------------\n$xml2 = @\"\n<w:p>\n  <w:pPr><w:pStyle w:val=\"Heading3\"/></w:pPr>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:r><w:t>Shayaan</w:t></w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n  <w:r><w:t xml:space=\"preserve\"> Ali, </w:t></w:r>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:r><w:t>Davinderpal</w:t></w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n  <w:r><w:t xml:space=\"preserve\"> Cheema, Alexander Powell, </w:t></w:r>\n  <w:proofErr w:type=\"spellStart\"/>\n  <w:r><w:t>Riddhi</w:t></w:r>\n  <w:proofErr w:type=\"spellEnd\"/>\n  <w:r><w:t xml:space=\"preserve\"> Sharma</w:t></w:r>\n</w:p>\n\"@\nSet-ParagraphOoxml 27 $xml2\n\n# --- 3. Introduction paragraph: wrap both \"tac\" occurrences (COM #31) ------\n$xml3 = @\"\n<w:p>\n  <w:r><w:tab/><w:t xml:space=\"preserve\">Tic Tac Toe is a game traditionally played between two players on a 3x3 grid. Players take turns placing X and O in the empty spaces on the grid. This project takes the basics of tic </w:t></w:r>\n  <w:proofErr w:type=\"gramStart\"/>\n  <w:r><w:t>tac</w:t></w:r>\n  <w:proofErr w:type=\"gramEnd\"/>\n  <w:r><w:t xml:space=\"preserve\"> toe and builds an application to allow players to play against each other over the internet. </w:t></w:r>\n  <w:r><w:t xml:space=\"preserve\">It will take the basics of tic </w:t></w:r>\n  <w:proofErr w:type=\"gramStart\"/>\n  <w:r><w:t>tac</w:t></w:r>\n  <w:proofErr w:type=\"gramEnd\"/>\n  <w:r><w:t xml:space=\"preserve\"> toe and stretch it to support any number of players on any board size (n x n). The application is built using a client-server architecture where clients connect to the server and play each other with the server acting as the messenger.</w:t></w:r>\n</w:p>\n\"@\nSet-ParagraphOoxml 31 $xml3\n\n# --- 4. Client Design paragraph: append a new sentence/run (COM #57) -------\n$xml4 = @\"\n<w:p>\n  <w:r><w:tab/><w:t>The client is designed and built to be a thin client, the client does not perform game calculations. The client is simply a user interface that accepts inputs from the user an converts them into game move. The game moves are then sent to the server where the server handles the game logic and sends board updates and game updates to the client.</w:t></w:r>\n  <w:r><w:t xml:space=\"preserve\"> The client uses a few threads, it has a thread for passing messages to the server and a thread that listens to the messages from the server as well as a thread for the GUI.</w:t></w:r>\n</w:p>\n\"@\nSet-ParagraphOoxml 57 $xml4\n\n# --- 5. Synchronization paragraph: split run, add a page break (COM #59) ---\n$xml5 = @\"\n<w:p>\n  <w:r><w:tab/><w:t xml:space=\"preserve\">The clients are synchronized by the server. The clients connect to the server and the server handles all message passing between them. If clients play out of turn then the server will stop them. If the clients try and make a move where another client has already moved then the client server will stop </w:t></w:r>\n  <w:r><w:lastRenderedPageBreak/><w:t>them and tell them to try again. The clients cannot become out of sync because all turn and game logic is computed on the server and the server facilitates all messages that are passed between the clients.</w:t></w:r>\n</w:p>\n\"@\nSet-ParagraphOoxml 59 $xml5\n\n# --- 6. \"Termination\" heading: drop the page break (it moved above) (#60) --\n$xml6 = @\"\n<w:p>\n  <w:pPr><w:pStyle w:val=\"Heading1\"/></w:pPr>\n  <w:r><w:t>Termination</w:t></w:r>\n</w:p>\n\"@\nSet-ParagraphOoxml 60 $xml6\n\n# --- 7. Termination detail paragraph: add the _GoBack bookmark (COM #61) ---\n$xml7 = @\"\n<w:p>\n  <w:r><w:tab/><w:t xml:space=\"preserve\">The server is terminated through the console. The server admin can simply type &#8216;stop&#8217; to stop the server from running. This terminates all threads that are running o</w:t></w:r>\n  <w:r><w:t>n the server and boots any clients from in progress games.</w:t></w:r>\n  <w:r><w:t xml:space=\"preserve\"> Termination in python is as simple as calling threads destructor in C++. Because threads cannot be killed externally the server requires a termination event object that is shared by all threads. When the event is signaled all threads terminate themselves.</w:t></w:r>\n  <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n  <w:bookmarkEnd w:id=\"0\"/>\n</w:p>\n\"@\nSet-ParagraphOoxml 61 $xml7\n\n# --- 8. \"Object oriented style\" paragraph: split + \"second\"->\"final\", and --\n#        add a new paragraph right after it with the prototyping narrative.\n#        The old trailing _GoBack bookmark is intentionally dropped here\n#        since it was relocated to the Termination paragraph (step 7).\n$xml8 = @\"\n<w:p>\n  <w:r><w:tab/><w:t xml:space=\"preserve\">The object oriented style of the </w:t></w:r>\n  <w:r><w:t xml:space=\"preserve\">final </w:t></w:r>\n  <w:r><w:t>version of the project took longer to make but it was a better result. The code of the final version is much more readable than the first and all of the classes perform fewer functions working together to make the game work instead of a few classes doing all of the work. The first version had a few points of failure for that would take the entire application with it if it failed. In the final version some threads can fail without taking down the server.</w:t></w:r>\n</w:p>\n<w:p>\n  <w:r><w:tab/><w:t>During development a throwaway prototyping model was used to create the final version. There were several iterations in order to create code that is easy to read, efficient and modular. The first iteration was not very object oriented and code acted more like several scripts running on different threads. The next version was more object oriented and focused on how the specific parts of the server and client would need to be built in an object oriented project. The next iteration took the knowledge from the first iteration and tried to extend the server and client classes to work together with each other. Then in the next iteration we focused on how the client and server would interact with the game and the user interface. Finally taking forward the knowledge of the last iteration the project was started again from the ground up in order to remove the mistakes from the prototype phase and finally create a product that met the goals of the final project.</w:t></w:r>\n</w:p>\n\"@\nSet-ParagraphOoxml 64 $xml8\n\n# --- 9. \"What did we think of the labs for this course?\" gets a page break -\n# Paragraph count grew by one (step 8 added a paragraph): originally #69, now #70.\n$xml9 = @\"\n<w:p>\n  <w:r><w:lastRenderedPageBreak/><w:t>What did we think of the labs for this course?</w:t></w:r>\n</w:p>\n\"@\nSet-ParagraphOoxml 70 $xml9\n"}
